# Apply scheduled market-data refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 194.75
$ws.Range("I11").Value = 194.75
$ws.Range("K11").Value = 194.75
$ws.Range("M11").Value = -54.75

$ws.Range("H74").Value = 6500
$ws.Range("I74").Value = 6500
$ws.Range("K74").Value = 6500
$ws.Range("M74").Value = -5564

$ws.Range("H77").Value = 6500
$ws.Range("I77").Value = 6500
$ws.Range("K77").Value = 32500
$ws.Range("M77").Value = -27820

$ws.Range("H88").Value = 2498.8
$ws.Range("I88").Value = 3000
$ws.Range("K88").Value = 3000
$ws.Range("M88").Value = -2594

$ws.Range("H91").Value = 2498.8
$ws.Range("I91").Value = 3000
$ws.Range("K91").Value = 3000
$ws.Range("M91").Value = -1596

$ws.Range("H96").Value = 1923.4286
$ws.Range("I96").Value = 2240.6667
$ws.Range("K96").Value = 6722.000100000001
$ws.Range("M96").Value = -5349.000100000001

$ws.Range("H97").Value = 1200
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3600
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -4592

$ws.Range("H98").Value = 2514.3125
$ws.Range("I98").Value = 2410.2222
$ws.Range("K98").Value = 2410.2222
$ws.Range("M98").Value = -912.2222000000002

$ws.Range("H99").Value = 794
$ws.Range("I99").Value = 352.8
$ws.Range("K99").Value = 1058.4
$ws.Range("M99").Value = 439.5999999999999

$ws.Range("H100").Value = 2029.6666
$ws.Range("I100").Value = 903
$ws.Range("K100").Value = 903
$ws.Range("M100").Value = -362

$ws.Range("H101").Value = 330
$ws.Range("I101").Value = 330
$ws.Range("K101").Value = 990
$ws.Range("M101").Value = 632

$ws.Range("H103").Value = 905.4167
$ws.Range("I103").Value = 884.74194
$ws.Range("K103").Value = 2654.22582
$ws.Range("M103").Value = -2068.22582

$ws.Range("H122").Value = 2514.3125
$ws.Range("I122").Value = 2410.2222
$ws.Range("K122").Value = 7230.6666
$ws.Range("M122").Value = -4780.6666

$ws.Range("H136").Value = 68934.5
$ws.Range("J136").Value = 68934.5
$ws.Range("L136").Value = 68934.5
$ws.Range("N136").Value = -79134.5

$ws.Range("H138").Value = 3382.3044
$ws.Range("J138").Value = 2934.8572
$ws.Range("L138").Value = 8804.571599999999
$ws.Range("N138").Value = -19084.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 398351.06
$ws.Range("J2").Value = 1770.5714
$ws.Range("L2").Value = 1770.5714
$ws.Range("N2").Value = -1996.5714

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""

$ws.Range("H61").Value = 2756.0715
$ws.Range("I61").Value = 1508.6
$ws.Range("J61").Value = 5874.75
$ws.Range("K61").Value = 1508.6
$ws.Range("L61").Value = 5874.75
$ws.Range("M61").Value = -1296.6
$ws.Range("N61").Value = -6298.75

$ws.Range("H63").Value = 4780.875
$ws.Range("J63").Value = 1582.6666
$ws.Range("L63").Value = 1582.6666
$ws.Range("N63").Value = -2954.6666

$ws.Range("H66").Value = 4780.875
$ws.Range("J66").Value = 1582.6666
$ws.Range("L66").Value = 7913.333000000001
$ws.Range("N66").Value = -14777.333

$ws.Range("H116").Value = 398351.06
$ws.Range("J116").Value = 1770.5714
$ws.Range("L116").Value = 1770.5714
$ws.Range("N116").Value = -6358.5714

$ws.Range("H132").Value = 1542.6923
$ws.Range("I132").Value = 1242.3636
$ws.Range("J132").Value = 3194.5
$ws.Range("K132").Value = 3727.0908
$ws.Range("L132").Value = 9583.5
$ws.Range("M132").Value = -1197.0908
$ws.Range("N132").Value = -14643.5

$ws.Range("H136").Value = 2756.0715
$ws.Range("I136").Value = 1508.6
$ws.Range("J136").Value = 5874.75
$ws.Range("K136").Value = 4525.799999999999
$ws.Range("L136").Value = 17624.25
$ws.Range("M136").Value = -1975.799999999999
$ws.Range("N136").Value = -22724.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 398351.06
$ws.Range("J3").Value = 1770.5714
$ws.Range("L3").Value = 1770.5714
$ws.Range("N3").Value = -1998.5714

$ws.Range("H46").Value = 30000
$ws.Range("J46").Value = 30000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30596

$ws.Range("H82").Value = 36125
$ws.Range("I82").Value = 20000
$ws.Range("J82").Value = 41500
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 41500
$ws.Range("M82").Value = -19617
$ws.Range("N82").Value = -42266

$ws.Range("H85").Value = 36125
$ws.Range("I85").Value = 20000
$ws.Range("J85").Value = 41500
$ws.Range("K85").Value = 20000
$ws.Range("L85").Value = 41500
$ws.Range("M85").Value = -18674
$ws.Range("N85").Value = -44152

$ws.Range("H99").Value = 1350.5
$ws.Range("J99").Value = 1702
$ws.Range("L99").Value = 1702
$ws.Range("N99").Value = -4698

$ws.Range("H107").Value = 547.2857
$ws.Range("I107").Value = 422.3
$ws.Range("K107").Value = 422.3
$ws.Range("M107").Value = 1497.7

$ws.Range("H134").Value = 8700.3125
$ws.Range("I134").Value = 9815.571
$ws.Range("J134").Value = 893.5
$ws.Range("K134").Value = 29446.713
$ws.Range("L134").Value = 2680.5
$ws.Range("M134").Value = -26911.713
$ws.Range("N134").Value = -7750.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2277.92
$ws.Range("I31").Value = 1955.2307
$ws.Range("J31").Value = 2627.5
$ws.Range("K31").Value = 1955.2307
$ws.Range("L31").Value = 2627.5
$ws.Range("M31").Value = -1660.2307
$ws.Range("N31").Value = -3217.5

$ws.Range("H34").Value = 2277.92
$ws.Range("I34").Value = 1955.2307
$ws.Range("J34").Value = 2627.5
$ws.Range("K34").Value = 1955.2307
$ws.Range("L34").Value = 2627.5
$ws.Range("M34").Value = -1753.2307
$ws.Range("N34").Value = -3031.5

$ws.Range("H58").Value = 2902276
$ws.Range("I58").Value = 8697771
$ws.Range("J58").Value = 4528.2
$ws.Range("K58").Value = 8697771
$ws.Range("L58").Value = 4528.2
$ws.Range("M58").Value = -8697568
$ws.Range("N58").Value = -4934.2

$ws.Range("H62").Value = 2122.0833
$ws.Range("I62").Value = 2267.8572
$ws.Range("J62").Value = 1918
$ws.Range("K62").Value = 2267.8572
$ws.Range("L62").Value = 1918
$ws.Range("M62").Value = -1643.8572
$ws.Range("N62").Value = -3166

$ws.Range("H65").Value = 2122.0833
$ws.Range("I65").Value = 2267.8572
$ws.Range("J65").Value = 1918
$ws.Range("K65").Value = 11339.286
$ws.Range("L65").Value = 9590
$ws.Range("M65").Value = -8219.286
$ws.Range("N65").Value = -15830

$ws.Range("H99").Value = 2843.7
$ws.Range("I99").Value = 1822.8334
$ws.Range("K99").Value = 1822.8334
$ws.Range("M99").Value = -324.8334

$ws.Range("H107").Value = 701.3333
$ws.Range("I107").Value = 469.5
$ws.Range("K107").Value = 469.5
$ws.Range("M107").Value = 1450.5

$ws.Range("H126").Value = 2843.7
$ws.Range("I126").Value = 1822.8334
$ws.Range("K126").Value = 5468.5002
$ws.Range("M126").Value = -2998.5002

$ws.Range("H136").Value = 2902276
$ws.Range("I136").Value = 8697771
$ws.Range("J136").Value = 4528.2
$ws.Range("K136").Value = 26093313
$ws.Range("L136").Value = 13584.6
$ws.Range("M136").Value = -26090763
$ws.Range("N136").Value = -18684.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 166675900
$ws.Range("I17").Value = 333333800
$ws.Range("J17").Value = 18001.334
$ws.Range("K17").Value = 1000001400
$ws.Range("L17").Value = 54004.00199999999
$ws.Range("M17").Value = -1000001231
$ws.Range("N17").Value = -54342.00199999999

$ws.Range("H81").Value = 2075
$ws.Range("J81").Value = 2666.6667
$ws.Range("L81").Value = 8000.000100000001
$ws.Range("N81").Value = -10246.0001

$ws.Range("H84").Value = 2075
$ws.Range("J84").Value = 2666.6667
$ws.Range("L84").Value = 24000.0003
$ws.Range("N84").Value = -35232.0003

$ws.Range("H131").Value = 27554.115
$ws.Range("J131").Value = 32433.955
$ws.Range("L131").Value = 97301.86500000001
$ws.Range("N131").Value = -107381.865

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 36639.5
$ws.Range("J127").Value = 36639.5
$ws.Range("L127").Value = 36639.5
$ws.Range("N127").Value = -46559.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5057.385
$ws.Range("J122").Value = 5379.2
$ws.Range("L122").Value = 16137.6
$ws.Range("N122").Value = -21037.6

$ws.Range("H136").Value = 4907.316
$ws.Range("I136").Value = 3963.1667
$ws.Range("J136").Value = 6525.857
$ws.Range("K136").Value = 11889.5001
$ws.Range("L136").Value = 19577.571
$ws.Range("M136").Value = -9339.500100000001
$ws.Range("N136").Value = -24677.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2044002.4
$ws.Range("J3").Value = 55003
$ws.Range("L3").Value = 55003
$ws.Range("N3").Value = -55231

$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = ""

$ws.Range("H107").Value = 1236.7
$ws.Range("I107").Value = 692.8
$ws.Range("K107").Value = 2078.4
$ws.Range("M107").Value = -158.3999999999996

$ws.Range("H136").Value = 10685043
$ws.Range("I136").Value = 17362398
$ws.Range("K136").Value = 52087194
$ws.Range("M136").Value = -52084644
